$d = $word.ActiveDocument

# --- Step 1: add a trailing space to the "Development on the project began..." paragraph ---
$d.Content.Find.Execute(
    "etc.) .",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "etc.) . ", 2) | Out-Null

# --- Step 2: locate the "Memory module" paragraph and rework its text + formatting ---
$memParaIndex = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Memory module has been implemented*") {
        $memParaIndex = $i
    }
    $i = $i + 1
}

$memPara = $d.Paragraphs($memParaIndex)
$memRange = $memPara.Range
# Re-writing the whole paragraph range (which includes the trailing manual line
# break + paragraph mark) drops the manual line break and keeps a single
# paragraph mark, matching the target.
$memRange.Text = "  - Memory module has been implemented and tested. Additional testing maybe required."
$d.Paragraphs($memParaIndex).Format.FirstLineIndent = 36

# --- Step 3: insert the five new bullet paragraphs right after it ---
$insertAt = $d.Paragraphs($memParaIndex).Range.End

$rsquo = [char]0x2019
$newLines = @(
    "  - master control and controls for cars and players were implemented",
    "  - testing is still required on these modules",
    "  - The source code for the game is completed in 3 days",
    "  - Testing began on each component of the game",
    ("  - The game can now input the lives, and number of cars of each type. Objects can now be displayed. They still can" + $rsquo + "t move due to a bug.")
)

$paraIdx = $memParaIndex + 1
foreach ($line in $newLines) {
    $r = $d.Range($insertAt, $insertAt)
    $r.InsertParagraphAfter()
    $nr = $d.Range($insertAt, $insertAt)
    $nr.Text = $line
    $d.Paragraphs($paraIdx).Format.FirstLineIndent = 36
    $insertAt = $insertAt + $line.Length + 1
    $paraIdx = $paraIdx + 1
}

# --- Step 4: restore the trailing whitespace + manual line break paragraph ---
# (this reuses the paragraph that used to directly follow the "Memory module"
# paragraph, which was already empty aside from its paragraph mark, so its
# "no first-line-indent" formatting is already correct)
$vtab = [char]11
$tailPara = $d.Paragraphs($paraIdx)
$tailRange = $tailPara.Range
$tailRange.Text = "            " + $vtab
